# Update the "Förändrad" (Changed) date column (C2:C9) from 2023-10-08 (45207)
# to 2023-10-09 (45208) for all 8 data rows, as per the commit's automatic update.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 9; $row++) {
    $cell = $ws.Cells.Item($row, 3)  # Column C
    $cell.Value = 45208
}
